$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.986.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.404.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.80%  "

$ws.Range("E4").Value = "  -0.34%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.639"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.78%  "

$ws.Range("E15").Value = "  +1.88%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.762.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.393.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.938.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.07%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000108"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.19%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "271.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.75%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "

$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.09%  "

$ws.Range("E31").Value = "  -2.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "36.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0925"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.38%  "

$ws.Range("E35").Value = "  +3.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.51%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0364"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.51%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.108"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.234"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.28%  "

$ws.Range("E44").Value = "  -0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +47.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "118.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.641.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.67%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.92%  "
